$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: AoA 22 deg, elevator control now SI (was dbx_v1_s50_r16_a25_bX..., NO)
$ws.Range("A19").Value = "dbx_v1_s50_r16_a22_b0_da0_deX_dr0"
$ws.Range("B19").Value = "SI"

# Row 20: AoA 25 deg, elevator control SI (was dbx_v1_s50_r16_a30_bX..., replaced)
$ws.Range("A20").Value = "dbx_v1_s50_r16_a25_b0_da0_deX_dr0"

# Row 21: new row - AoA 35 deg, elevator control SI
$ws.Range("A21").Value = "dbx_v1_s50_r16_a35_b0_da0_deX_dr0"
$ws.Range("A21").WrapText = $true
$ws.Range("B21").Value = "SI"
$ws.Range("B21").WrapText = $true

# Apply underline Arial 10pt font to the whole column B data range -> new style (fontId 5)
$bRange = $ws.Range("B2:B21")
$bRange.Font.Underline = $true
$bRange.Font.Size = 10
$bRange.Font.Name = "Arial"

# Update selection to reflect the newly added data range
$ws.Range("B2:B21").Select()
